$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = $new
        Write-Output "OK: $old"
    } else {
        Write-Output "NOT FOUND: $old"
    }
}

# 1 & 8: Title text appears twice (Heading1 title, and bold run near the end) - both get the same new text
Replace-Text "Play Jimi Hendrix Slot for Free - Review of Game Theme, Gameplay and Features" "Play Jimi Hendrix for Free - Slot Game Review"
Replace-Text "Play Jimi Hendrix Slot for Free - Review of Game Theme, Gameplay and Features" "Play Jimi Hendrix for Free - Slot Game Review"

# 2: "What we like" bullet 1
Replace-Text "Captures the spirit of 60s revolution with design and soundtrack" "Captures the spirit of the 60s and Jimi Hendrix's iconic guitar style"

# 3: "What we like" bullet 2
Replace-Text "Personalized betting with adjustable bet levels and coin costs" "Thematic symbols and soundtrack create an immersive experience"

# 4: "What we like" bullet 3
Replace-Text "Six special features for extra excitement and chance to win big" "Multiple special features and bonus rounds for added excitement"

# 5: "What we like" bullet 4
Replace-Text "Works seamlessly on all devices with HTML5 technology" "Compatible with desktops, mobile phones, and tablets"

# 6: "What we don't like" bullet 1
Replace-Text "20 fixed pay lines limit betting flexibility" "Limited number of pay lines compared to some other slot games"

# 7: "What we don't like" bullet 2
Replace-Text "Only medium volatility, so might not appeal to higher risk players" "Some players may not be familiar with Jimi Hendrix's music"

# 9: Italic summary paragraph at end
Replace-Text "Experience the spirit of 60s with Jimi Hendrix, a 5-reel slot with 20 fixed pay lines. Play for free and enjoy six special features on all devices." "Read our review of Jimi Hendrix, a slot game that captures the spirit of the 60s. Play for free to experience the iconic guitar style and special features."
